# Auto-generated edit script: updates market-data derived columns (H-N)
# across the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed scheduled pull of Leve profit data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 20838792
$ws.Range("J33").Value = 6571.7144
$ws.Range("L33").Value = 6571.7144
$ws.Range("N33").Value = -7029.7144

$ws.Range("H43").Value = 1001
$ws.Range("I43").Value = 1001
$ws.Range("K43").Value = 1001
$ws.Range("M43").Value = -932

$ws.Range("H116").Value = 5948.857
$ws.Range("I116").Value = 6333.3335
$ws.Range("J116").Value = 5844
$ws.Range("K116").Value = 6333.3335
$ws.Range("L116").Value = 5844
$ws.Range("M116").Value = -2891.3335
$ws.Range("N116").Value = -12728

$ws.Range("H135").Value = 1835.2941
$ws.Range("I135").Value = 1835.2941
$ws.Range("K135").Value = 16517.6469
$ws.Range("M135").Value = -13982.6469

$ws.Range("H138").Value = 1649.41
$ws.Range("I138").Value = 1102
$ws.Range("J138").Value = 2079.5178
$ws.Range("K138").Value = 3306
$ws.Range("L138").Value = 6238.553400000001
$ws.Range("M138").Value = 1834
$ws.Range("N138").Value = -16518.5534

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3137534.5
$ws.Range("I110").Value = 3369884.8
$ws.Range("K110").Value = 3369884.8
$ws.Range("M110").Value = -3367839.8

$ws.Range("H132").Value = 2953416.5
$ws.Range("I132").Value = 3224.76
$ws.Range("K132").Value = 9674.280000000001
$ws.Range("M132").Value = -7144.280000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 26646.555
$ws.Range("I99").Value = 30842.6
$ws.Range("J99").Value = 5666.3335
$ws.Range("K99").Value = 30842.6
$ws.Range("L99").Value = 5666.3335
$ws.Range("M99").Value = -29344.6
$ws.Range("N99").Value = -8662.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 16753
$ws.Range("I86").Value = 19327.715
$ws.Range("K86").Value = 19327.715
$ws.Range("M86").Value = -18204.715

$ws.Range("H89").Value = 16753
$ws.Range("I89").Value = 19327.715
$ws.Range("K89").Value = 96638.575
$ws.Range("M89").Value = -91022.575

$ws.Range("H94").Value = 1047.7778
$ws.Range("I94").Value = 861
$ws.Range("J94").Value = 1101.1428
$ws.Range("K94").Value = 861
$ws.Range("L94").Value = 1101.1428
$ws.Range("M94").Value = -410
$ws.Range("N94").Value = -2003.1428

$ws.Range("H97").Value = 16800
$ws.Range("J97").Value = 16800
$ws.Range("L97").Value = 16800
$ws.Range("N97").Value = -18782

$ws.Range("H107").Value = 1516.4166
$ws.Range("I107").Value = 2157
$ws.Range("J107").Value = 619.6
$ws.Range("K107").Value = 2157
$ws.Range("L107").Value = 619.6
$ws.Range("M107").Value = -237
$ws.Range("N107").Value = -4459.6

$ws.Range("H121").Value = 100000
$ws.Range("J121").Value = 100000
$ws.Range("L121").Value = 100000
$ws.Range("N121").Value = -102620

$ws.Range("H132").Value = 43481216
$ws.Range("I132").Value = 2847.9
$ws.Range("K132").Value = 8543.700000000001
$ws.Range("M132").Value = -6013.700000000001

$ws.Range("H134").Value = 45463544
$ws.Range("I134").Value = 2408.625
$ws.Range("J134").Value = 71441336
$ws.Range("K134").Value = 7225.875
$ws.Range("L134").Value = 214324008
$ws.Range("M134").Value = -4690.875
$ws.Range("N134").Value = -214329078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1616059.9
$ws.Range("J4").Value = 666828
$ws.Range("L4").Value = 2000484
$ws.Range("N4").Value = -2000708

$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 6000
$ws.Range("M64").Value = -5730

$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 6000
$ws.Range("M67").Value = -5064

$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2189

$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4944

$ws.Range("H80").Value = 11416.667
$ws.Range("I80").Value = 5750
$ws.Range("J80").Value = 14250
$ws.Range("K80").Value = 17250
$ws.Range("L80").Value = 42750
$ws.Range("M80").Value = -16314
$ws.Range("N80").Value = -44622

$ws.Range("H83").Value = 11416.667
$ws.Range("I83").Value = 5750
$ws.Range("J83").Value = 14250
$ws.Range("K83").Value = 51750
$ws.Range("L83").Value = 128250
$ws.Range("M83").Value = -47070
$ws.Range("N83").Value = -137610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 667.875
$ws.Range("I107").Value = 524.36365
$ws.Range("K107").Value = 524.36365
$ws.Range("M107").Value = 1395.63635

$ws.Range("H122").Value = 1789720.9
$ws.Range("I122").Value = 1789720.9
$ws.Range("K122").Value = 5369162.699999999
$ws.Range("M122").Value = -5366712.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4717
$ws.Range("I61").Value = 3485.875
$ws.Range("K61").Value = 3485.875
$ws.Range("M61").Value = -3283.875

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H113").Value = 4717
$ws.Range("I113").Value = 3485.875
$ws.Range("K113").Value = 3485.875
$ws.Range("M113").Value = -1315.875

$ws.Range("H136").Value = 13256.314
$ws.Range("I136").Value = 10080.64
$ws.Range("K136").Value = 30241.92
$ws.Range("M136").Value = -27691.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 717.86957
$ws.Range("I107").Value = 732.36365
$ws.Range("J107").Value = 399
$ws.Range("K107").Value = 2197.09095
$ws.Range("L107").Value = 1197
$ws.Range("M107").Value = -277.0909499999998
$ws.Range("N107").Value = -5037

$ws.Range("H122").Value = 604507.0600000001
$ws.Range("I122").Value = 857229.6
$ws.Range("J122").Value = 7162.727
$ws.Range("K122").Value = 2571688.8
$ws.Range("L122").Value = 21488.181
$ws.Range("M122").Value = -2569238.8
$ws.Range("N122").Value = -26388.181

$ws.Range("H132").Value = 14013.523
$ws.Range("I132").Value = 6750.64
$ws.Range("J132").Value = 24694.234
$ws.Range("K132").Value = 20251.92
$ws.Range("L132").Value = 74082.702
$ws.Range("M132").Value = -17721.92
$ws.Range("N132").Value = -79142.702

$ws.Range("H136").Value = 18126.215
$ws.Range("J136").Value = 30692.066
$ws.Range("L136").Value = 92076.198
$ws.Range("N136").Value = -97176.198
